$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update rows 2-7: merge the two alternating speed conditions into a
# single value (1.5) for columns A (nr) and B (nm), and rewrite the
# remaining per-row figures so the table now represents 4 loops. ---

# Row 2
$ws.Range("A2").Value = 1.5
$ws.Range("B2").Value = 1.5

# Row 3
$ws.Range("A3").Value = 1.5
$ws.Range("B3").Value = 1.5
$ws.Range("C3").Value = 3
$ws.Range("E3").Value = 2

# Row 4
$ws.Range("A4").Value = 1.5
$ws.Range("B4").Value = 1.5
$ws.Range("C4").Value = 99
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 2.5

# Row 5
$ws.Range("A5").Value = 1.5
$ws.Range("B5").Value = 1.5
$ws.Range("C5").Value = 99
$ws.Range("E5").Value = 6

# Row 6
$ws.Range("A6").Value = 1.5
$ws.Range("B6").Value = 1.5
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 1

# Row 7
$ws.Range("A7").Value = 1.5
$ws.Range("B7").Value = 1.5
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 15

# --- Remove the now-unneeded extra loop rows (8-13); their content is
# cleared, leaving the lower blank/styled rows (21-23) untouched. ---
$ws.Range("A8:G13").ClearContents() | Out-Null

# --- Update the active selection to match where the user left off. ---
$ws.Range("F8").Select()
